$d = $word.ActiveDocument

# 1. Re-create the "smarthosting" bookmark so it is (re)assigned a fresh,
#    low bookmark id (matches w:id="0" in the target document instead of
#    the original w:id="1").
$bm = $d.Bookmarks.Item("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# 2. Remove the two HYPERLINK field codes (CryptoBridge and HitBTC) that
#    sit inside the "100,000 Smart:" bullet paragraph. Deleting the Field
#    objects removes the begin/separate/end fldChar runs, the hidden
#    instrText run and the visible link-text run in one shot, leaving no
#    empty shell runs behind. Delete the second one first so the index of
#    the first one is unaffected.
$fields = $d.Fields
$fields.Item(2).Delete()
$d.Fields.Item(1).Delete()

# 3. Clean up the now-orphaned surrounding text: " such as" (incl. the
#    trailing non-breaking space) and the bold ", " + non-breaking-space
#    run that used to separate the two links collapse away, leaving
#    "...obtained from exchanges" directly followed by the existing bold
#    "." run.
$cleanupRange = $d.Paragraphs(13).Range
$found = $cleanupRange.Find.Execute(
    " such as" + [char]0x00A0 + "," + [char]0x00A0,
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $cleanupRange.Delete()
}
